# Add "Unlock Spells" section (3 spells + 4 status strings) to the strings
# workbook: 10 new rows (139-148) on Sheet1, plus a comment on A139.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row data -----------------------------------------------------
# (StringId, Text) pairs, in the order they appear starting at row 139 /
# A-value 6610137.
$newRows = @(
    "Arcane Unlock",
    "A short burst of arcane magic will deal with most simple locks.",
    "Finesse Unlock",
    "More stubborn locks require some finesse. This spell is able to manipulate some of the inner components allowing opening of more complex locks.",
    "Power Unlock",
    "This spell allows you to channel more of your power into the lock, opening more complex varieties.",
    "Spell too weak to unlock",
    "Lock cannot be unlocked with magic",
    "Wrong type of target to unlock",
    "Target already unlocked"
)

$startRow = 139
$startId = 6610137

# --- Apply formatting first by copying from the existing rows so the new
# cells pick up the same style indices used by the rest of the sheet
# (A column: fill used by rows like A3:A5; B column: wrap-text style used
# throughout the list, e.g. B138).
$lastRow = $startRow + $newRows.Count - 1

$ws.Range("A3").Copy()
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null

$ws.Range("B138").Copy()
$ws.Range("B" + $startRow + ":B" + $lastRow).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Fill in the values -------------------------------------------------
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startId + $i
    $ws.Cells.Item($row, 2).Value = $newRows[$i]
}

# --- Add the cell comment on A139 ---------------------------------------
$excel.UserName = "Jim"
$ws.Range("A139").AddComment("Unlock Spells") | Out-Null

Write-Output "Added rows $startRow to $lastRow and comment on A139"
